$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "63.821.22"
Set-TextValue "E2" "  +1.09%  "
Set-TextValue "D3" "3.208.89"
Set-TextValue "E3" "  -2.58%  "
Set-TextValue "E4" "  -0.18%  "
Set-TextValue "D5" "595.49"
Set-TextValue "E5" "  -1.33%  "
Set-TextValue "D6" "137.42"
Set-TextValue "E6" "  -1.53%  "
Set-TextValue "E7" "  -0.06%  "
Set-TextValue "D8" "3.206.85"
Set-TextValue "E8" "  -2.53%  "
Set-TextValue "D9" "0.516"
Set-TextValue "E9" "  +0.30%  "
Set-TextValue "E10" "  -3.25%  "
Set-TextValue "D11" "5.29"
Set-TextValue "E11" "  -3.58%  "
Set-TextValue "D12" "0.458"
Set-TextValue "E12" "  -1.28%  "
Set-TextValue "D13" "0.0000242"
Set-TextValue "E13" "  -0.82%  "
Set-TextValue "D14" "34.97"
Set-TextValue "E14" "  +1.74%  "
Set-TextValue "D15" "3.739.17"
Set-TextValue "E15" "  -2.58%  "
Set-TextValue "E16" "  -1.52%  "
Set-TextValue "D17" "3.197.16"
Set-TextValue "E17" "  -3.26%  "
Set-TextValue "D18" "63.665.24"
Set-TextValue "E18" "  +0.61%  "
Set-TextValue "D19" "6.63"
Set-TextValue "E19" "  -2.38%  "
Set-TextValue "D20" "467.39"
Set-TextValue "E20" "  -1.63%  "
Set-TextValue "D21" "14.07"
Set-TextValue "E21" "  +1.38%  "
Set-TextValue "D22" "0.705"
Set-TextValue "E22" "  -3.46%  "
Set-TextValue "D23" "7.70"
Set-TextValue "E23" "  -2.61%  "
Set-TextValue "D24" "13.58"
Set-TextValue "E24" "  -1.47%  "
Set-TextValue "D25" "83.19"
Set-TextValue "E25" "  -2.09%  "
Set-TextValue "E26" "  +0.31%  "
Set-TextValue "D27" "2.70"
Set-TextValue "E27" "  -2.11%  "
Set-TextValue "D28" "0.996"
Set-TextValue "E28" "  -0.37%  "
Set-TextValue "D29" "7.85"
Set-TextValue "E29" "  -2.59%  "
Set-TextValue "D30" "6.87"
Set-TextValue "E30" "  -3.42%  "
Set-TextValue "D31" "2.07"
Set-TextValue "E31" "  -2.40%  "
Set-TextValue "D32" "27.56"
Set-TextValue "E32" "  -2.92%  "
Set-TextValue "D33" "0.104"
Set-TextValue "E33" "  -0.65%  "
Set-TextValue "D34" "2.42"
Set-TextValue "E34" "  -2.71%  "
Set-TextValue "E35" "  -4.35%  "
Set-TextValue "D36" "5.89"
Set-TextValue "E36" "  -1.48%  "
Set-TextValue "D37" "51.66"
Set-TextValue "E37" "  -0.82%  "
Set-TextValue "D38" "0.0₃0732"
Set-TextValue "E38" "  +0.46%  "
Set-TextValue "D39" "0.0392"
Set-TextValue "E39" "  -1.38%  "
Set-TextValue "D40" "2.75"
Set-TextValue "E40" "  +2.93%  "
Set-TextValue "D41" "408.78"
Set-TextValue "E41" "  -3.81%  "
Set-TextValue "D42" "8.15"
Set-TextValue "E42" "  -1.24%  "
Set-TextValue "D43" "0.113"
Set-TextValue "E43" "  -4.69%  "
Set-TextValue "D44" "2.831.26"
Set-TextValue "E44" "  -9.43%  "
Set-TextValue "D45" "0.257"
Set-TextValue "E45" "  -1.27%  "
Set-TextValue "D46" "2.17"
Set-TextValue "E46" "  -0.97%  "
Set-TextValue "B47" "Arweave"
Set-TextValue "C47" "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D47" "35.98"
Set-TextValue "E47" "  -0.74%  "
Set-TextValue "B48" "USDe"
Set-TextValue "C48" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D48" "0.999"
Set-TextValue "E48" "  +0.01%  "
Set-TextValue "B49" "Monero"
Set-TextValue "C49" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D49" "127.06"
Set-TextValue "E49" "  -0.55%  "
Set-TextValue "D50" "25.72"
Set-TextValue "E50" "  -1.35%  "
Set-TextValue "E51" "  -0.60%  "
